# Applies:
#  1. Rename header "Requested quantity" -> "Weekly_PO_Qty" on "Weekly Quantity" sheet (B1)
#  2. Rename header "Requested quantity" -> "Monthly_PO_Qty" on "Monthly Trend" sheet (B1)
#  3. Add a new worksheet "PO Forecast" after "Monthly Trend" with forecast data

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# 1) Update header text on the Weekly Quantity sheet
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# 2) Update header text on the Monthly Trend sheet
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 3) Add the new "PO Forecast" sheet as the last sheet (after Monthly Trend)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the page margins used on the other sheets in this workbook
$wsForecast.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$wsForecast.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$wsForecast.PageSetup.TopMargin = $excel.InchesToPoints(1)
$wsForecast.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$wsForecast.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$wsForecast.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)

# Copy header formatting (bold, centered, bordered) from the Weekly Quantity header row
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Copy date-cell formatting from the Weekly Quantity date column onto column A of the new sheet
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A61").PasteSpecial(-4122)  # xlPasteFormats

# Header row values
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
    @(44997.99999999999, 23, -15.61506340866121, 57.09786530621413),
    @(45011.99999999999, 24, -11.76873756880311, 60.65042121521493),
    @(45018.99999999999, 24, -10.12525676165679, 60.56215488677437),
    @(45025.99999999999, 24, -9.83851055539499, 58.75412967137404),
    @(45032.99999999999, 24, -9.974103823685299, 58.63932091176819),
    @(45039.99999999999, 25, -14.2183116898207, 59.07501646975803),
    @(45046.99999999999, 25, -14.94945315942833, 57.54133688259731),
    @(45053.99999999999, 25, -11.23073210037863, 57.03637928588709),
    @(45060.99999999999, 25, -10.8179939662141, 62.06929623279016),
    @(45067.99999999999, 25, -10.66370749177719, 60.56027882229612),
    @(45074.99999999999, 26, -9.493105806318853, 62.31552837046255),
    @(45088.99999999999, 26, -11.82441277675334, 60.7557079536303),
    @(45116.99999999999, 27, -7.066349811587548, 64.73607840145843),
    @(45123.99999999999, 27, -9.852328678775018, 64.28308017955324),
    @(45144.99999999999, 28, -6.019080940277938, 62.6479663483631),
    @(45165.99999999999, 29, -6.51066407129093, 65.55852719797404),
    @(45172.99999999999, 29, -7.014186971367607, 64.50367467456202),
    @(45186.99999999999, 29, -7.911254452358056, 66.4555029913241),
    @(45193.99999999999, 30, -6.042749756182386, 66.19386444557448),
    @(45200.99999999999, 30, -6.055707797051016, 63.7291700121112),
    @(45214.99999999999, 30, -5.192378473933331, 64.2916366606765),
    @(45221.99999999999, 31, -6.189659976969986, 66.69307911957254),
    @(45235.99999999999, 31, -6.508018660634038, 64.22381174451584),
    @(45242.99999999999, 31, -5.123927747894798, 66.38673495779123),
    @(45256.99999999999, 32, -3.866557331032655, 67.16130378604058),
    @(45263.99999999999, 32, -2.830672303494972, 67.39284117914821),
    @(45270.99999999999, 32, -2.776873682707441, 69.76145306631385),
    @(45312.99999999999, 34, -2.609990160351677, 68.76874221418564),
    @(45319.99999999999, 34, -2.556384995684348, 70.64444196615726),
    @(45340.99999999999, 34, -1.196992733213019, 70.76076133019515),
    @(45347.99999999999, 35, -3.302979800750293, 69.67551181021553),
    @(45354.99999999999, 35, -1.442816884516019, 68.8431203335129),
    @(45361.99999999999, 35, -0.0737346713457338, 72.87636836793122),
    @(45389.99999999999, 36, 0.6635189858648661, 70.40980776563792),
    @(45403.99999999999, 37, -0.06948297716695957, 73.68723740350283),
    @(45410.99999999999, 37, 2.516840273167112, 74.01715897332009),
    @(45417.99999999999, 37, 2.220476379602284, 72.4752500235456),
    @(45445.99999999999, 38, 0.921064154097934, 72.93177178013423),
    @(45452.99999999999, 38, 2.661541520604036, 73.51230760743267),
    @(45459.99999999999, 38, 1.858980583165666, 75.07369457658973),
    @(45466.99999999999, 39, 3.29763374330522, 74.19138823773828),
    @(45487.99999999999, 39, 4.561100818299553, 75.75252711077322),
    @(45501.99999999999, 40, 4.204437177600406, 73.94866140182279),
    @(45515.99999999999, 40, 2.039836089043189, 75.13217048735099),
    @(45536.99999999999, 41, 7.330910265056937, 75.93426628069044),
    @(45550.99999999999, 41, 7.786195265500148, 76.6929271506157),
    @(45557.99999999999, 42, 5.102379065313959, 76.4975782569702),
    @(45564.99999999999, 42, 5.391485907680945, 76.49075532776756),
    @(45578.99999999999, 42, 8.508292369086112, 77.81215929688048),
    @(45585.99999999999, 43, 6.404679038589513, 80.3956891781029),
    @(45592.99999999999, 43, 6.787046990523588, 78.09803923762061),
    @(45599.99999999999, 43, 7.329516356577764, 76.75273496161992),
    @(45606.99999999999, 43, 7.192023774475824, 79.22643348599517),
    @(45613.99999999999, 43, 8.54868367854486, 77.56201054777954),
    @(45620.99999999999, 44, 8.790613507961309, 79.94224567319669),
    @(45627.99999999999, 44, 7.873272221309033, 79.3230277862793),
    @(45634.99999999999, 44, 10.17405603984611, 81.86657864494735),
    @(45641.99999999999, 44, 8.948007748334366, 78.79868699898168),
    @(45648.99999999999, 45, 7.940001570890627, 80.56460792244373),
    @(45655.99999999999, 45, 12.51110260152766, 78.98169426803778)
)

$row = 2
foreach ($entry in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $entry[0]
    $wsForecast.Cells.Item($row, 2).Value = $entry[1]
    $wsForecast.Cells.Item($row, 3).Value = $entry[2]
    $wsForecast.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}

Write-Host "Applied PO Forecast sheet and header renames."
